# benchmarks/template.xlsx -- "results with 4.3.2, fixed ft reporting"
#
# The underlying edit is a data refresh: new simulation results were pasted
# into C9:C23 (raw per-event times), the event count (C5) and the ms/event
# reference (E5) were updated, and the subsystem row labels for rows 16-21
# were re-ordered/fixed (this is the "fixed ft reporting" part of the
# commit message). All downstream formulas (D:L) are left untouched and
# simply recompute. Finally F5:G5 (the grand totals) get a "0.0" number
# format on top of their existing highlighted "Note" cell style, and the
# active selection moves to D32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Number of events / ms-per-event reference -----------------------------
$ws.Range("C5").Value = 10000
$ws.Range("E5").Value = 10000

# --- Raw per-subsystem times pasted from the new results.txt ---------------
$ws.Range("C9").Value  = 11.8891
$ws.Range("C10").Value = 20.1439
$ws.Range("C11").Value = 24.0514
$ws.Range("C12").Value = 27.7781
$ws.Range("C13").Value = 354.774
$ws.Range("C14").Value = 585.203
$ws.Range("C15").Value = 807.851
$ws.Range("C16").Value = 1805.18
$ws.Range("C17").Value = 1995.62
$ws.Range("C18").Value = 2249.58
$ws.Range("C19").Value = 2315.2
$ws.Range("C20").Value = 2315.2
$ws.Range("C21").Value = 2798.42
$ws.Range("C22").Value = 4188.15
$ws.Range("C23").Value = 5752.89

# --- Subsystem labels for rows 16-21 re-ordered / fixed ("ft" reporting) ---
$ws.Range("B16").Value = "torus"
$ws.Range("B17").Value = "rich"
$ws.Range("B18").Value = "ltcc"
$ws.Range("B19").Value = "ftof"
$ws.Range("B20").Value = "ft"
$ws.Range("B21").Value = "dc"

# --- Highlight the grand totals with a one-decimal number format -----------
$ws.Range("F5:G5").NumberFormat = "0.0"

# --- Move the active selection (matches the saved view state) --------------
$ws.Range("D32").Select()
